# Trace_Report_SBHP_initial.xlsx edit
# Commit message: "cahged white color code if" (sic) - underlying data update:
#  - refreshed the search-completed timestamp in the report header (A1)
#  - swapped rows 4 and 5 data: car CGEX1941 (Junction Delivery to JOHNSTOWN via GWR,
#    GREELEY, CO) is now row 4, and car CAIX541012 (Departure to LOVELAND via HKCKDE,
#    LA JUNTA, CO) is now row 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / search timestamp text
$ws.Range("A1").Value = "Description unknown, completed 06/22/2023 11:08:20 EDT, by WPJTOWN1.The search returned: 4 events."

# Row 4 -> now the CGEX1941 car record
$ws.Range("A4").Value = "CGEX"
$ws.Range("B4").Value = 1941
$ws.Range("C4").Value = "GREELEY"
$ws.Range("D4").Value = "CO"
$ws.Range("G4").Value = 1419
$ws.Range("H4").Value = "Junction Delivery"
$ws.Range("I4").Value = "GWR"
$ws.Range("J4").Value = "JOHNSTOWN"
$ws.Range("L4").Value = 198750
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 198750
$ws.Range("O4").Value = "CGEX1941"

# Row 5 -> now the CAIX541012 car record
$ws.Range("A5").Value = "CAIX"
$ws.Range("B5").Value = 541012
$ws.Range("C5").Value = "LA JUNTA"
$ws.Range("F5").Value = 22
$ws.Range("G5").Value = 719
$ws.Range("H5").Value = "Departure"
$ws.Range("I5").Value = "HKCKDE"
$ws.Range("J5").Value = "LOVELAND"
$ws.Range("L5").Value = 273100
$ws.Range("M5").Value = 62900
$ws.Range("N5").Value = 210200
$ws.Range("O5").Value = "CAIX541012"
